# Insert a new data row at row 59 (pushes existing rows 59..144 down to 60..145,
# matching the dimension change from A1:R144 to A1:R145), then populate the
# newly inserted row with the new record's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 59, shifting the rest of the table down.
$ws.Rows.Item(59).Insert()

# Populate the new row 59 with the new record.
$ws.Range("A59").Value = 10
$ws.Range("B59").Value = "Vega Modelo de Temuco"
$ws.Range("C59").Value = "La Araucanía"
$ws.Range("D59").Value = 44494
$ws.Range("E59").Value = 9
$ws.Range("F59").Value = 100112013
$ws.Range("G59").Value = "Alcachofa"
$ws.Range("H59").Value = "Española"
$ws.Range("I59").Value = "Extra"
$ws.Range("J59").Value = 1000
$ws.Range("K59").Value = 500
$ws.Range("L59").Value = 500
$ws.Range("M59").Value = 500
$ws.Range("N59").Value = "$/unidad"
$ws.Range("O59").Value = "Región Metropolitana"
$ws.Range("P59").Value = 500
$ws.Range("Q59").Value = 1
$ws.Range("R59").Value = "Hortaliza"
